$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I19").Value = -0.9853713629864949
$ws.Range("J19").Value = 0.2062375548920074
$ws.Range("K19").Value = 0.2052426328312739
$ws.Range("L19").Value = 2.39656107586889
